# Auto-generated: update crypto price table with latest values scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="26.541.47"; E="  +0.79%  "}
    @{Row=3; D="1.812.92"; E="  +1.11%  "}
    @{Row=4; D="1.006"; E="  -0.14%  "}
    @{Row=5; E="  -0.04%  "}
    @{Row=6; D="305.23"; E="  -0.47%  "}
    @{Row=7; D="0.4644"; E="  +2.11%  "}
    @{Row=8; D="0.3575"; E="  -1.25%  "}
    @{Row=9; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.07100"; E="  +0.44%  "}
    @{Row=10; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.8981"; E="  +3.13%  "}
    @{Row=11; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.07765"; E="  -0.22%  "}
    @{Row=12; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="19.34"; E="  -0.03%  "}
    @{Row=13; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.859.22"; E="  +4.68%  "}
    @{Row=14; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="5.232"; E="  -0.53%  "}
    @{Row=15; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="6.301"; E="  -0.19%  "}
    @{Row=16; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="87.46"; E="  +3.37%  "}
    @{Row=17; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.008"; E="  -0.01%  "}
    @{Row=18; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.000008532"; E="  +0.41%  "}
    @{Row=19; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.006"; E="  -0.10%  "}
    @{Row=20; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="26.587.29"; E="  +0.84%  "}
    @{Row=21; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="14.15"; E="  +0.10%  "}
    @{Row=22; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="4.983"; E="  +0.19%  "}
    @{Row=23; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="10.53"; E="  +0.49%  "}
    @{Row=24; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="1.917"; E="  -2.96%  "}
    @{Row=25; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="151.79"; E="  -0.07%  "}
    @{Row=26; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="17.84"; E="  +0.24%  "}
    @{Row=27; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.985"; E="  -2.96%  "}
    @{Row=28; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="112.86"; E="  +0.47%  "}
    @{Row=29; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="4.797"; E="  -0.65%  "}
    @{Row=30; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.08728"; E="  +0.89%  "}
    @{Row=31; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="3.122"; E="  +3.12%  "}
    @{Row=32; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="2.759"; E="  +2.96%  "}
    @{Row=33; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.7306"; E="  +2.32%  "}
    @{Row=34; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.420"; E="  -0.43%  "}
    @{Row=35; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.119"; E="  +1.13%  "}
    @{Row=36; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.074"; E="  -0.38%  "}
    @{Row=37; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01927"; E="  -0.68%  "}
    @{Row=38; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.915"; E="  +1.82%  "}
    @{Row=39; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.05090"; E="  +0.13%  "}
    @{Row=40; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.5045"; E="  +2.58%  "}
    @{Row=41; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="6.775"; E="  -1.50%  "}
    @{Row=42; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1491"; E="  -1.44%  "}
    @{Row=43; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="7.952"; E="  -0.21%  "}
    @{Row=44; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.4669"; E="  +2.19%  "}
    @{Row=45; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="1.006"; E="  -0.03%  "}
    @{Row=46; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="9.997"; E="  +1.84%  "}
    @{Row=47; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="98.23"; E="  -1.24%  "}
    @{Row=48; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="1.562"; E="  -1.23%  "}
    @{Row=49; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.06006"; E="  +1.01%  "}
    @{Row=50; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="63.55"; E="  +0.15%  "}
    @{Row=51; B="Elrond"; C="https://coinranking.com/coin/omwkOTglq+elrond-egld"; D="35.71"; E="  -0.90%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Force the price column to remain plain text (matches source file, which never
        # stores these as numeric cells) by using a leading quote-prefix, then stripping
        # the style Excel auto-applies for quote-prefixed text so formatting stays default.
        $ws.Cells.Item($r, 4).Value = "'" + $u.D
        $ws.Cells.Item($r, 4).ClearFormats()
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
